$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.164.47'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '3.578.46'
$ws.Range("E3").Value = '  -1.93%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.36%  '
$ws.Range("D7").Value = '3.576.52'
$ws.Range("E7").Value = '  -1.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.619'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.36%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.184'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.652'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000306'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.32%  '
$ws.Range("D15").Value = '4.153.58'
$ws.Range("E15").Value = '  -1.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.12%  '
$ws.Range("D17").Value = '3.575.70'
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("D18").Value = '70.027.58'
$ws.Range("E18").Value = '  -1.31%  '
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("E21").Value = '  -3.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '489.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '95.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.61%  '
$ws.Range("E28").Value = '  -6.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.76'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '66.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.44%  '
$ws.Range("E34").Value = '  -6.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '571.95'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +15.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.419'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '38.75'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.46%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").Value = '0.0₃0793'
$ws.Range("E40").Value = '  -4.67%  '
$ws.Range("E41").Value = '  -3.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("D45").Value = '3.218.64'
$ws.Range("E45").Value = '  -3.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0443'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.59'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +32.12%  '
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
